# Generate Report for Handback
# Updates timestamps + priority recorded during a later handback run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-10-21 01:29:02"
$wsOverview.Range("G3").Value = "2016-10-21 01:29:02"

# zh-cn sheet: Priority + Handoff/Handback datetimes
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-10-21 01:28:50"
$wsZhCn.Range("H3").Value = "2016-10-21 01:28:50"
$wsZhCn.Range("K2").Value = "2016-10-21 01:29:34"
$wsZhCn.Range("K3").Value = "2016-10-21 01:29:34"

# de-de sheet: Priority + Handoff Generate Date + Handback datetime
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-10-21 01:29:02"
$wsDeDe.Range("H3").Value = "2016-10-21 01:29:02"
$wsDeDe.Range("K2").Value = "2016-10-21 01:29:52"
$wsDeDe.Range("K3").Value = "2016-10-21 01:29:52"
